$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 52: 2025-10-06 data point.
# Force column A to be stored as text (not auto-parsed into a date serial)
# by briefly applying a text number format, then restore the default
# "Normal" style so no extra style index is left on the cell.
$ws.Range("A52").NumberFormat = "@"
$ws.Range("A52").Value = "2025-10-06"
$ws.Range("A52").Style = "Normal"

$ws.Range("B52").Value = 54.11000061035156
$ws.Range("C52").Value = 712.6500244140625
$ws.Range("D52").Value = 335.1000061035156
